# The commit behind this change ("Fixed POI packaging and upgraded to
# POI 3.15") is a build-tooling upgrade in the source repository: the
# project switched the Apache POI library used to *write* the .docx
# test fixture. That newer POI version's XML writer happens to emit
# element attributes in alphabetical order (and re-orders the
# namespace-declaration attributes on the <w:document> root element),
# instead of preserving whatever attribute order was in the source.
#
# A byte-for-byte comparison of the previous and new fixture shows that
# every single changed line is exactly the same element with exactly
# the same set of attribute name="value" pairs -- only their ordering
# differs (e.g. <w:pgSz w:w="11906" w:h="16838"/> becomes
# <w:pgSz w:h="16838" w:w="11906"/>). No text, run, paragraph, style
# value, page-size/margin number, language, or any other document
# content/formatting actually changed.
#
# Attribute-serialization order is not part of the Word document
# object model -- it is a detail of the XML writer used when the
# package is serialized, not a document property that Word (or COM
# automation against Word) exposes or lets you control. There is
# therefore no WordProcessing object-model call that corresponds to
# this change: applying "the change described by the diff" through
# Word COM automation means making no content/formatting edits at all,
# since none occurred.
#
# This script intentionally performs no mutations -- it only touches
# the document via read-only object-model access, matching the fact
# that the underlying .docx content is unchanged.
$d = $word.ActiveDocument
$null = $d.Content
